$wb = $excel.ActiveWorkbook

# --- Change 1: "sets" sheet, D4 (home_points for set 3) 5 -> 6 ---
$wsSets = $wb.Worksheets.Item("sets")
$wsSets.Range("D4").Value = 6

# --- Change 2: "rallies" sheet, append new row 69 ---
$wsRallies = $wb.Worksheets.Item("rallies")
$wsRallies.Range("A69").Value = 68
$wsRallies.Range("B69").Value = 1
$wsRallies.Range("C69").Value = 3
$wsRallies.Range("D69").Value = 6
$wsRallies.Range("E69").Value = "NOS"
$wsRallies.Range("G69").Value = 4
$wsRallies.Range("H69").Value = "MEIO"
$wsRallies.Range("I69").Value = "PONTO"
$wsRallies.Range("J69").Value = "NOS"
$wsRallies.Range("K69").Value = 6
$wsRallies.Range("L69").Value = 0
$wsRallies.Range("M69").Value = "1 4 m"
$wsRallies.Range("N69").Value = "FRENTE"
$wsRallies.Range("O69").Value = "FRENTE"
$wsRallies.Range("P69").Value = "FRENTE"
